{"js": "// Dairy levy template update: the merge-field placeholder\n// \"{d.DairyTestDataLoadDate}\" is being renamed to \"{d.PreviousMonth}\".\n// Find the run containing the old field name and swap in the new one,\n// leaving every other run / property untouched.\nconst results = context.document.body.search(\"DairyTestDataLoadDate\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Expected to find \"DairyTestDataLoadDate\" in the document body, but found none.');\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"PreviousMonth\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Dairy levy template update: the merge-field placeholder\n# \"{d.DairyTestDataLoadDate}\" is being renamed to \"{d.PreviousMonth}\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"DairyTestDataLoadDate\"\n$find.Replacement.Text = \"PreviousMonth\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n)\n"}
